$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (single-dot decimal-looking strings),
# so they match the original inline-string formatting.
$textCells = @("D5", "D6", "D9", "D11", "D14", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D31", "D33", "D34", "D35", "D36", "D39", "D41", "D42", "D46", "D50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "61.787.54"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "3.416.25"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "578.17"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").Value = "145.03"
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "7.65"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").Value = "0.386"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "3.999.14"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "28.15"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "3.427.28"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").Value = "0.0000170"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "61.808.00"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "6.16"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "13.89"
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").Value = "9.18"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("D21").Value = "390.31"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").Value = "74.27"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").Value = "0.551"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "0.0000115"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").Value = "0.188"
$ws.Range("E26").Value = "  +2.81%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "7.46"
$ws.Range("E27").Value = "  +3.23%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "8.00"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").Value = "23.49"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").Value = "5.23"
$ws.Range("E34").Value = "  +5.16%  "
$ws.Range("D35").Value = "6.96"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "168.62"
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("D37").Value = "3.448.46"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").Value = "28.60"
$ws.Range("E39").Value = "  +6.55%  "
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D41").Value = "0.787"
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("D42").Value = "4.45"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("E44").Value = "  +4.39%  "
$ws.Range("D45").Value = "2.507.97"
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("D46").Value = "22.77"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "2.11"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("E51").Value = "  -0.06%  "
